$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N ("Late"), shifting the
# existing N/O/P columns (Late / Outstanding-heading / Outstanding) one
# place to the right -> O/P/Q.
$ws.Columns("N").Insert()

# The newly inserted column inherits the width Excel would naturally give
# it (same width as the column to its left, but without the "best fit"
# flag baked in).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and select cell K15, as
# in the authored workbook.
$ws.Activate()
$ws.Range("K15").Select()
